$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2025-09-20 Saturday", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-09-21 Sunday", 2)

# Update the answer table. Cells are addressed by (row, column) rather than by
# their old text, because one value ("83÷8=10, 3") occurs twice in the table
# and each occurrence maps to a different new value.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "44÷5=8, 4"
$t.Cell(1, 2).Range.Text  = "16÷8=2, 0"
$t.Cell(1, 3).Range.Text  = "77÷3=25, 2"
$t.Cell(1, 4).Range.Text  = "86÷8=10, 6"
$t.Cell(1, 5).Range.Text  = "24÷4=6, 0"

$t.Cell(5, 1).Range.Text  = "41÷9=4, 5"
$t.Cell(5, 2).Range.Text  = "14÷5=2, 4"
$t.Cell(5, 3).Range.Text  = "64÷5=12, 4"
$t.Cell(5, 4).Range.Text  = "57÷3=19, 0"
$t.Cell(5, 5).Range.Text  = "40÷3=13, 1"

$t.Cell(9, 1).Range.Text  = "11÷4=2, 3"
$t.Cell(9, 2).Range.Text  = "93÷9=10, 3"
$t.Cell(9, 3).Range.Text  = "34÷8=4, 2"
$t.Cell(9, 4).Range.Text  = "50÷8=6, 2"
$t.Cell(9, 5).Range.Text  = "71÷4=17, 3"

$t.Cell(13, 1).Range.Text = "68÷2=34, 0"
$t.Cell(13, 2).Range.Text = "38÷4=9, 2"
$t.Cell(13, 3).Range.Text = "33÷8=4, 1"
$t.Cell(13, 4).Range.Text = "60÷4=15, 0"
$t.Cell(13, 5).Range.Text = "83÷4=20, 3"

$t.Cell(17, 1).Range.Text = "53÷4=13, 1"
$t.Cell(17, 2).Range.Text = "48÷4=12, 0"
$t.Cell(17, 3).Range.Text = "38÷8=4, 6"
$t.Cell(17, 4).Range.Text = "37÷9=4, 1"
$t.Cell(17, 5).Range.Text = "73÷3=24, 1"
